$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price cells we are about to update keep their text representation
# even though many of the new values look numeric (Excel would otherwise
# auto-convert them to real numbers). D4 and D21 are left untouched since
# their Price value does not change.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7:D20").NumberFormat = "@"
$ws.Range("D22:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.431.72'
$ws.Range("E2").Value = '  -0.21%  '

$ws.Range("D3").Value = '2.101.62'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '335.24'
$ws.Range("E5").Value = '  +1.59%  '

$ws.Range("D7").Value = '0.5230'
$ws.Range("E7").Value = '  -0.46%  '

$ws.Range("D8").Value = '0.4567'
$ws.Range("E8").Value = '  +4.05%  '

$ws.Range("D9").Value = '56.44'
$ws.Range("E9").Value = '  +14.66%  '

$ws.Range("D10").Value = '0.08938'
$ws.Range("E10").Value = '  +0.63%  '

$ws.Range("D11").Value = '1.179'
$ws.Range("E11").Value = '  +1.19%  '

$ws.Range("D12").Value = '24.21'
$ws.Range("E12").Value = '  -2.69%  '

$ws.Range("D13").Value = '2.097.19'
$ws.Range("E13").Value = '  -0.22%  '

$ws.Range("D14").Value = '6.845'
$ws.Range("E14").Value = '  +1.48%  '

$ws.Range("D15").Value = '8.055'
$ws.Range("E15").Value = '  +3.74%  '

$ws.Range("D16").Value = '97.43'
$ws.Range("E16").Value = '  +0.86%  '

$ws.Range("D17").Value = '0.00001153'
$ws.Range("E17").Value = '  +1.97%  '

$ws.Range("D18").Value = '1.005'
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").Value = '0.06641'
$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").Value = '19.21'
$ws.Range("E20").Value = '  -0.49%  '

$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").Value = '6.308'
$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("D23").Value = '30.495.69'
$ws.Range("E23").Value = '  -0.21%  '

$ws.Range("D24").Value = '12.40'
$ws.Range("E24").Value = '  +1.02%  '

$ws.Range("D25").Value = '2.354'
$ws.Range("E25").Value = '  +0.29%  '

$ws.Range("D26").Value = '2.342.95'
$ws.Range("E26").Value = '  -0.25%  '

$ws.Range("D27").Value = '22.21'
$ws.Range("E27").Value = '  -1.21%  '

$ws.Range("D28").Value = '162.94'
$ws.Range("E28").Value = '  +0.48%  '

$ws.Range("D29").Value = '2.520'
$ws.Range("E29").Value = '  -4.31%  '

$ws.Range("D30").Value = '133.46'
$ws.Range("E30").Value = '  +0.30%  '

$ws.Range("D31").Value = '1.211'
$ws.Range("E31").Value = '  -1.13%  '

$ws.Range("D32").Value = '0.1070'
$ws.Range("E32").Value = '  -0.31%  '

$ws.Range("D33").Value = '1.657'
$ws.Range("E33").Value = '  -1.97%  '

$ws.Range("D34").Value = '6.364'
$ws.Range("E34").Value = '  +1.94%  '

$ws.Range("D35").Value = '3.947'
$ws.Range("E35").Value = '  +1.40%  '

$ws.Range("D36").Value = '10.31'
$ws.Range("E36").Value = '  +0.97%  '

$ws.Range("D37").Value = '5.944'
$ws.Range("E37").Value = '  +7.66%  '

$ws.Range("D38").Value = '0.02582'
$ws.Range("E38").Value = '  -0.36%  '

$ws.Range("D39").Value = '0.06862'
$ws.Range("E39").Value = '  +1.74%  '

$ws.Range("D40").Value = '0.2338'
$ws.Range("E40").Value = '  +2.23%  '

$ws.Range("D41").Value = '12.67'
$ws.Range("E41").Value = '  -0.81%  '

$ws.Range("D42").Value = '0.6881'
$ws.Range("E42").Value = '  -0.75%  '

$ws.Range("D43").Value = '1.248'
$ws.Range("E43").Value = '  -2.15%  '

$ws.Range("D44").Value = '2.331'
$ws.Range("E44").Value = '  +4.64%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '14.07'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6404'
$ws.Range("E46").Value = '  -0.49%  '

$ws.Range("D47").Value = '3.661'
$ws.Range("E47").Value = '  +0.68%  '

$ws.Range("D48").Value = '1.249'
$ws.Range("E48").Value = '  -0.30%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.00000000342'
$ws.Range("E49").Value = '  +14.76%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '83.13'
$ws.Range("E50").Value = '  +0.30%  '

$ws.Range("D51").Value = '1.203'
$ws.Range("E51").Value = '  -1.23%  '

